$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) updates: force Text format so numeric-looking
# strings (with significant trailing/leading zeros) are preserved exactly
# as authored, matching the source inlineStr cells.
$priceCells = @{
    "D2"  = "243.48"
    "D3"  = "23.85"
    "D4"  = "5.242"
    "D5"  = "0.05871"
    "D6"  = "6.464"
    "D7"  = "3.334"
    "D9"  = "0.8776"
    "D10" = "0.1380"
    "D11" = "0.07257"
    "D12" = "0.03070"
    "D13" = "0.03053"
    "D14" = "0.09323"
    "D15" = "3.851"
    "D16" = "0.001542"
    "D17" = "0.04696"
    "D18" = "0.0006001"
    "D19" = "0.006273"
    "D20" = "0.001263"
    "D21" = "0.004578"
    "D22" = "0.00008701"
    "D24" = "2.177"
    "D25" = "0.3205"
    "D28" = "0.0002340"
    "D40" = "0.03787"
    "D41" = "0.006321"
    "D42" = "0.1050"
    "D44" = "0.007787"
    "D45" = "0.00005495"
    "D47" = "0.5401"
    "D48" = "0.02176"
}

foreach ($addr in $priceCells.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceCells[$addr]
}

# --- Rows 10-18: the coin ranking list shifted by one position (One
# dropped from rank 9 down to rank 17), so Coin/Link/Volume columns need
# to be rewritten to reflect the new row contents.
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("E10").Value = "9WazirXWRX"

$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"

$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"

$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("E13").Value = "12BitrueCoinBTR"

$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("E14").Value = "13BitMartTokenBMX"

$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("E15").Value = "14MCDexMCB"

$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("E16").Value = "15BitForexTokenBF"

$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("E17").Value = "16CoinExTokenCET"

$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("E18").Value = "17OneONE"

# --- Misc "Worst/Best in 24h" badge text churn on the Volume column.
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"
